# Fixed naive component forecaster bug - Presentation state 11.02.
# A new data point was inserted at the top of the results table (row 2),
# pushing all existing rows down by one and dropping the oldest row
# (former row 11) off the bottom of the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture current (pre-edit) values for columns B:G, rows 2-11, in order.
$oldValues = @()
for ($r = 2; $r -le 11; $r++) {
    $row = @()
    for ($c = 2; $c -le 7; $c++) {
        $row += $ws.Cells.Item($r, $c).Value()
    }
    $oldValues += ,$row
}

# Shift every row's B:G values down into the next row (row 11's old data
# falls off the bottom of the table and is discarded).
for ($r = 11; $r -ge 3; $r--) {
    $srcRow = $oldValues[$r - 3]
    for ($c = 2; $c -le 7; $c++) {
        $ws.Cells.Item($r, $c).Value = $srcRow[$c - 2]
    }
}

# Write the brand-new data point into row 2.
$ws.Range("B2").Value = -0.02843597334849807
$ws.Range("C2").Value = 0.5120978290967556
$ws.Range("D2").Value = 0.6508982077648369
$ws.Range("E2").Value = 0.8067826273320694
$ws.Range("F2").Value = 0.8283753329562445
$ws.Range("G2").Value = 19
